$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (header in row 1 stays the same).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics.
$data = @(
    @("ECs",  "Cx3cl1", "Itgav", "ECs",  2, 0.6666666666666666, 6.662215666666666,  19.986647,  0.4964100230945138, 0.4964100230945138, 3, 1, 16.535604,          49.606812,   0.2120453146491552, 0.2120453146491552, 110.163760026596,  991.4738402393639, 0.1052614195420706, 0.1052614195420706),
    @("ECs",  "Cx3cl1", "Itgav", "FAPs", 2, 0.6666666666666666, 6.662215666666666,  19.986647,  0.4964100230945138, 0.4964100230945138, 3, 1, 40.62063066666667, 121.861892,  0.5209011059384622, 0.5209011059384622, 270.6234020173471, 2435.610618156124, 0.2585805300288698, 0.2585805300288698),
    @("ECs",  "Cx3cl1", "Itgav", "sCs",  2, 0.6666666666666666, 6.662215666666666,  19.986647,  0.4964100230945138, 0.4964100230945138, 3, 1, 20.825229,         62.475687,   0.2670535794123827, 0.2670535794123827, 138.742166905721,  1248.679502151489, 0.1325680735235735, 0.1325680735235735),
    @("FAPs", "Cx3cl1", "Itgav", "ECs",  3, 1,                  5.098534333333333, 15.295603,  0.3798981709375522, 0.3798981709375522, 3, 1, 16.535604,          49.606812,   0.2120453146491552, 0.2120453146491552, 84.30734471640399, 758.766102447636,  0.08055562719109181, 0.08055562719109183),
    @("FAPs", "Cx3cl1", "Itgav", "FAPs", 3, 1,                  5.098534333333333, 15.295603,  0.3798981709375522, 0.3798981709375522, 3, 1, 40.62063066666667, 121.861892,  0.5209011059384622, 0.5209011059384622, 207.1056800956529, 1863.951120860876, 0.1978893773853699, 0.1978893773853699),
    @("FAPs", "Cx3cl1", "Itgav", "sCs",  3, 1,                  5.098534333333333, 15.295603,  0.3798981709375522, 0.3798981709375522, 3, 1, 20.825229,         62.475687,   0.2670535794123827, 0.2670535794123827, 106.178145056029,  955.603305504261,  0.1014531663610905, 0.1014531663610905),
    @("sCs",  "Cx3cl1", "Itgav", "ECs",  3, 1,                  1.660042,          4.980126,   0.1236918059679339, 0.123691805967934,  3, 1, 16.535604,          49.606812,   0.2120453146491552, 0.2120453146491552, 27.449797135368,   247.048174218312,  0.02622826791599281, 0.02622826791599281),
    @("sCs",  "Cx3cl1", "Itgav", "FAPs", 3, 1,                  1.660042,          4.980126,   0.1236918059679339, 0.123691805967934,  3, 1, 40.62063066666667, 121.861892,  0.5209011059384622, 0.5209011059384622, 67.43195297315468, 606.8875767583921, 0.06443119852422247, 0.06443119852422248),
    @("sCs",  "Cx3cl1", "Itgav", "sCs",  3, 1,                  1.660042,          4.980126,   0.1236918059679339, 0.123691805967934,  3, 1, 20.825229,         62.475687,   0.2670535794123827, 0.2670535794123827, 34.570754799618,   311.136793196562,  0.03303233952771868, 0.03303233952771868)
)

$startRow = 2
$numCols = $data[0].Length
for ($c = 0; $c -lt $numCols; $c++) {
    for ($i = 0; $i -lt $data.Length; $i++) {
        $rowIndex = $startRow + $i
        $ws.Cells.Item($rowIndex, $c + 1).Value = $data[$i][$c]
    }
}

$ws.Range("A1:T10").Columns.AutoFit() | Out-Null
